$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for a new row ("Monday t+1") at row 3 -----------------------
# Shift existing rows 3..16 down to 4..17 by copying whole rows (values +
# formats) starting from the bottom so we never overwrite data we still need.
for ($r = 16; $r -ge 3; $r--) {
    $src = $ws.Rows.Item($r)
    $dst = $ws.Rows.Item($r + 1)
    $src.Copy()
    $dst.PasteSpecial(-4104)
}
$ws.Application.CutCopyMode = $false

# Re-apply the bold/bordered label formatting (same as the rest of column A)
# to the freshly vacated row 3 (and make sure the whole label column stays
# consistent after the row shuffle above).
$ws.Range("A2").Copy()
$ws.Range("A3:A17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Write the final labels / statistics for every data row ----------------
$ws.Range("A2").Value2 = "intercept"
$ws.Range("B2").Value2 = -0.02085163863032587
$ws.Range("C2").Value2 = 0.0004041259086361501
$ws.Range("D2").Value2 = -51.59688647702959
$ws.Range("E2").Value2 = 0

$ws.Range("A3").Value2 = "Monday t+1"
$ws.Range("B3").Value2 = 0.001275778340500457
$ws.Range("C3").Value2 = 0.001281899858242169
$ws.Range("D3").Value2 = 0.995224652142402
$ws.Range("E3").Value2 = 0.3196685464228193

$ws.Range("A4").Value2 = "SMB_loading"
$ws.Range("B4").Value2 = 0.0002553380462654783
$ws.Range("C4").Value2 = 0.0000827131283235668
$ws.Range("D4").Value2 = 3.087031665234778
$ws.Range("E4").Value2 = 0.00203120536490777

$ws.Range("A5").Value2 = "HML_loading"
$ws.Range("B5").Value2 = -0.00006496217218491791
$ws.Range("C5").Value2 = 0.00006719639272818322
$ws.Range("D5").Value2 = -0.9667508856866324
$ws.Range("E5").Value2 = 0.3337089044612326

$ws.Range("A6").Value2 = "RMW_loading"
$ws.Range("B6").Value2 = -0.00002613883641151181
$ws.Range("C6").Value2 = 0.00008071462826042827
$ws.Range("D6").Value2 = -0.3238426165722282
$ws.Range("E6").Value2 = 0.7460688705531968

$ws.Range("A7").Value2 = "CMA_loading"
$ws.Range("B7").Value2 = 0.000008921368959778128
$ws.Range("C7").Value2 = 0.00005444234911318493
$ws.Range("D7").Value2 = 0.1638681854309908
$ws.Range("E7").Value2 = 0.8698406245290997

$ws.Range("A8").Value2 = "momentum_loading"
$ws.Range("B8").Value2 = 0.0002206321355144677
$ws.Range("C8").Value2 = 0.0004387759296475268
$ws.Range("D8").Value2 = 0.5028355490961042
$ws.Range("E8").Value2 = 0.6150989966988687

$ws.Range("A9").Value2 = "reversal_loading"
$ws.Range("B9").Value2 = -0.0003298379531944742
$ws.Range("C9").Value2 = 0.0001092008857139268
$ws.Range("D9").Value2 = -3.020469578044903
$ws.Range("E9").Value2 = 0.002534831234876874

$ws.Range("A10").Value2 = "mkt_loading"
$ws.Range("B10").Value2 = 0.0005027548708982777
$ws.Range("C10").Value2 = 0.0001459215719862755
$ws.Range("D10").Value2 = 3.445377294493263
$ws.Range("E10").Value2 = 0.0005743060261175503

$ws.Range("A11").Value2 = "size*Monday_loading"
$ws.Range("B11").Value2 = -0.00003294250670592486
$ws.Range("C11").Value2 = 0.000133312957115638
$ws.Range("D11").Value2 = -0.2471065635229286
$ws.Range("E11").Value2 = 0.804834503250093

$ws.Range("A12").Value2 = "BM*Monday_loading"
$ws.Range("B12").Value2 = 0.00002165216185534479
$ws.Range("C12").Value2 = 0.0001076587937680762
$ws.Range("D12").Value2 = 0.2011183768414583
$ws.Range("E12").Value2 = 0.8406130676762107

$ws.Range("A13").Value2 = "ROE*Monday_loading"
$ws.Range("B13").Value2 = 0.000006368560262832938
$ws.Range("C13").Value2 = 0.0000373641148200928
$ws.Range("D13").Value2 = 0.17044590226471
$ws.Range("E13").Value2 = 0.8646654225861503

$ws.Range("A14").Value2 = "INV*Monday_loading"
$ws.Range("B14").Value2 = 0.00001023433002935715
$ws.Range("C14").Value2 = 0.00001152262201096253
$ws.Range("D14").Value2 = 0.8881945463124877
$ws.Range("E14").Value2 = 0.3744729562457861

$ws.Range("A15").Value2 = "MOM*Monday_loading"
$ws.Range("B15").Value2 = 0.001191056157142948
$ws.Range("C15").Value2 = 0.0006312761368349262
$ws.Range("D15").Value2 = 1.886743514675892
$ws.Range("E15").Value2 = 0.05924475211405898

$ws.Range("A16").Value2 = "REV*Monday_loading"
$ws.Range("B16").Value2 = 0.0003881395485549126
$ws.Range("C16").Value2 = 0.0002490685060889527
$ws.Range("D16").Value2 = 1.558364622849152
$ws.Range("E16").Value2 = 0.1192014162802011

$ws.Range("A17").Value2 = "mkt_risk_premium*Monday_loading"
$ws.Range("B17").Value2 = -0.00002404073062694617
$ws.Range("C17").Value2 = 0.0002379400461058029
$ws.Range("D17").Value2 = -0.1010369251431353
$ws.Range("E17").Value2 = 0.919524641667423
